# removed the cfg file from pdfgen
#
# The "pdfgen" worksheet had a row describing a "compare" step against
# /web/admin/pdf.cfg (with an ignore-list value of "warn"). That row is
# removed entirely, and the remaining "compare" row's Ignore List value
# is changed from "warn" to "report".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("pdfgen")
$ws.Activate()

# Update the Ignore List value for the remaining "compare" row (row 6)
$ws.Range("D6").Value = "report"

# Remove the whole row for the /web/admin/pdf.cfg entry (row 7), shifting
# the rows below it up.
$ws.Rows("7:7").Delete()

# Restore the selection to match the saved worksheet state.
$ws.Range("B17").Select()
